$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "21×41=861"
$t.Rows.Item(1).Cells.Item(2).Range.Text = "73×26=1898"
$t.Rows.Item(1).Cells.Item(3).Range.Text = "70×72=5040"
$t.Rows.Item(1).Cells.Item(4).Range.Text = "42×81=3402"
$t.Rows.Item(1).Cells.Item(5).Range.Text = "91×66=6006"
$t.Rows.Item(5).Cells.Item(1).Range.Text = "54×83=4482"
$t.Rows.Item(5).Cells.Item(2).Range.Text = "79×11=869"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "52×78=4056"
$t.Rows.Item(5).Cells.Item(4).Range.Text = "83×84=6972"
$t.Rows.Item(5).Cells.Item(5).Range.Text = "29×68=1972"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "26×38=988"
$t.Rows.Item(10).Cells.Item(2).Range.Text = "51×67=3417"
$t.Rows.Item(10).Cells.Item(3).Range.Text = "62×96=5952"
$t.Rows.Item(10).Cells.Item(4).Range.Text = "88×99=8712"
$t.Rows.Item(10).Cells.Item(5).Range.Text = "36×61=2196"
$t.Rows.Item(15).Cells.Item(1).Range.Text = "15×45=675"
$t.Rows.Item(15).Cells.Item(2).Range.Text = "79×16=1264"
$t.Rows.Item(15).Cells.Item(3).Range.Text = "22×94=2068"
$t.Rows.Item(15).Cells.Item(4).Range.Text = "78×25=1950"
$t.Rows.Item(15).Cells.Item(5).Range.Text = "11×20=220"
$t.Rows.Item(20).Cells.Item(1).Range.Text = "33×97=3201"
$t.Rows.Item(20).Cells.Item(2).Range.Text = "92×64=5888"
$t.Rows.Item(20).Cells.Item(3).Range.Text = "13×42=546"
$t.Rows.Item(20).Cells.Item(4).Range.Text = "74×11=814"
$t.Rows.Item(20).Cells.Item(5).Range.Text = "32×99=3168"
